$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header / data in column D (Tn) ---
$ws.Range("D2").Value = "Tn"
$ws.Range("D3").Value = 6
$ws.Range("D4").Value = 6
$ws.Range("D5").Value = 6
$ws.Range("D6").Value = 6
$ws.Range("D8").Value = 6
$ws.Range("D9").Value = 6
$ws.Range("D10").Value = 6

# Row 9 column H was missing, now has a value, which lets the shared
# formula in I9 compute instead of raising #DIV/0!
$ws.Range("H9").Value = 2716.682

# --- New columns J / K (duplicate timing headers for "Fixed buffer" block) ---
$ws.Range("J1").Value = "Fixed buffer"
$ws.Range("J2").Value = "Kernel time (ms)"
$ws.Range("K2").Value = "GOP/s"

$ws.Range("J3").Value = 7807.82
$ws.Range("J4").Value = 5519.646
$ws.Range("J5").Value = 4905.685
$ws.Range("J9").Value = 4749.926

# --- Formatting ---
# J1 gets a yellow fill (new style)
$ws.Range("J1").Interior.Color = 65535

# D9 switches from the old "center" style to a "right" aligned style
$ws.Range("D9").HorizontalAlignment = -4152  # xlRight

# --- Sheet view: scroll so column B is the leftmost visible column, and
#     move the active selection to J6 ---
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("J6").Select()
